{"js": "// Replace each two-digit multiplication problem/answer with its new value.\n// The document is a single table of \"AA\u00d7BB=CCCC\" style text runs; every\n// occurrence is unique, so a literal, case-sensitive search-and-replace\n// for each pair is unambiguous and safe.\nconst replacements = [\n  [\"86\u00d754=4644\", \"73\u00d780=5840\"],\n  [\"68\u00d716=1088\", \"69\u00d761=4209\"],\n  [\"11\u00d714=154\", \"26\u00d721=546\"],\n  [\"76\u00d716=1216\", \"21\u00d742=882\"],\n  [\"67\u00d774=4958\", \"51\u00d741=2091\"],\n  [\"22\u00d740=880\", \"85\u00d766=5610\"],\n  [\"80\u00d733=2640\", \"83\u00d742=3486\"],\n  [\"12\u00d721=252\", \"31\u00d761=1891\"],\n  [\"55\u00d799=5445\", \"92\u00d750=4600\"],\n  [\"51\u00d735=1785\", \"86\u00d753=4558\"],\n  [\"44\u00d779=3476\", \"26\u00d732=832\"],\n  [\"70\u00d737=2590\", \"46\u00d787=4002\"],\n  [\"54\u00d789=4806\", \"57\u00d768=3876\"],\n  [\"99\u00d789=8811\", \"63\u00d720=1260\"],\n  [\"24\u00d762=1488\", \"71\u00d775=5325\"],\n  [\"47\u00d795=4465\", \"11\u00d736=396\"],\n  [\"87\u00d724=2088\", \"27\u00d738=1026\"],\n  [\"26\u00d712=312\", \"93\u00d770=6510\"],\n  [\"41\u00d735=1435\", \"62\u00d764=3968\"],\n  [\"97\u00d776=7372\", \"23\u00d774=1702\"],\n  [\"26\u00d717=442\", \"27\u00d740=1080\"],\n  [\"63\u00d745=2835\", \"85\u00d724=2040\"],\n  [\"79\u00d718=1422\", \"56\u00d767=3752\"],\n  [\"20\u00d783=1660\", \"20\u00d731=620\"],\n  [\"79\u00d711=869\", \"30\u00d767=2010\"],\n];\n\nconst body = context.document.body;\nconst searchResults = [];\n\nfor (const [oldText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  searchResults.push(found);\n}\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, newText] = replacements[i];\n  const found = searchResults[i];\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Replace each two-digit multiplication problem/answer with its new value.\n# The document is a single table of \"AA\u00d7BB=CCCC\" style text runs; every\n# occurrence is unique, so a literal, case-sensitive Find/Replace for each\n# pair is unambiguous and safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"86\u00d754=4644\", \"73\u00d780=5840\"),\n    @(\"68\u00d716=1088\", \"69\u00d761=4209\"),\n    @(\"11\u00d714=154\", \"26\u00d721=546\"),\n    @(\"76\u00d716=1216\", \"21\u00d742=882\"),\n    @(\"67\u00d774=4958\", \"51\u00d741=2091\"),\n    @(\"22\u00d740=880\", \"85\u00d766=5610\"),\n    @(\"80\u00d733=2640\", \"83\u00d742=3486\"),\n    @(\"12\u00d721=252\", \"31\u00d761=1891\"),\n    @(\"55\u00d799=5445\", \"92\u00d750=4600\"),\n    @(\"51\u00d735=1785\", \"86\u00d753=4558\"),\n    @(\"44\u00d779=3476\", \"26\u00d732=832\"),\n    @(\"70\u00d737=2590\", \"46\u00d787=4002\"),\n    @(\"54\u00d789=4806\", \"57\u00d768=3876\"),\n    @(\"99\u00d789=8811\", \"63\u00d720=1260\"),\n    @(\"24\u00d762=1488\", \"71\u00d775=5325\"),\n    @(\"47\u00d795=4465\", \"11\u00d736=396\"),\n    @(\"87\u00d724=2088\", \"27\u00d738=1026\"),\n    @(\"26\u00d712=312\", \"93\u00d770=6510\"),\n    @(\"41\u00d735=1435\", \"62\u00d764=3968\"),\n    @(\"97\u00d776=7372\", \"23\u00d774=1702\"),\n    @(\"26\u00d717=442\", \"27\u00d740=1080\"),\n    @(\"63\u00d745=2835\", \"85\u00d724=2040\"),\n    @(\"79\u00d718=1422\", \"56\u00d767=3752\"),\n    @(\"20\u00d783=1660\", \"20\u00d731=620\"),\n    @(\"79\u00d711=869\", \"30\u00d767=2010\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
